{"js": "// Remove the trailing \"Ver no Jupiter...\" and \"\u00a9 2020 ...\" paragraphs\n// (plus the blank paragraph that precedes them), which sat right after\n// the \"LOQ4064: ...\" requirements line, leaving that line followed\n// directly by the blank paragraph / page break that used to trail them.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"LOQ4064\" requirements paragraph.\nlet reqIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOQ4064\") !== -1) {\n    reqIndex = i;\n    break;\n  }\n}\n\nif (reqIndex === -1) {\n  throw new Error(\"Could not find the 'LOQ4064' requirements paragraph.\");\n}\n\n// The three paragraphs right after it are the blank spacer, the\n// \"Ver no Jupiter...\" line and the \"\u00a9 2020 ...\" footer line. Delete them.\nconst toDelete = [];\nfor (let i = reqIndex + 1; i < items.length; i++) {\n  const text = items[i].text;\n  if (\n    text.trim() === \"\" ||\n    text.indexOf(\"Ver no Jupiter\") !== -1 ||\n    text.indexOf(\"Powered by Jekyll\") !== -1\n  ) {\n    toDelete.push(items[i]);\n    if (text.indexOf(\"Powered by Jekyll\") !== -1) {\n      break;\n    }\n  } else {\n    break;\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" and \"(c) 2020 ...\" paragraphs\n# (plus the blank paragraph that precedes them), which sat right after\n# the \"LOQ4064: ...\" requirements line, leaving that line followed\n# directly by the blank paragraph / page break that used to trail them.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n\n# Locate the \"LOQ4064\" requirements paragraph.\n$reqIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -match \"LOQ4064\") {\n        $reqIndex = $i\n        break\n    }\n}\n\nif ($reqIndex -eq -1) {\n    throw \"Could not find the 'LOQ4064' requirements paragraph.\"\n}\n\n# Find the last paragraph of the trailing block: the footer line\n# containing \"Powered by Jekyll\" (the \"(c) 2020 ...\" paragraph).\n$endIndex = -1\nfor ($i = $reqIndex + 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -match \"Powered by Jekyll\") {\n        $endIndex = $i\n        break\n    }\n}\n\nif ($endIndex -eq -1) {\n    throw \"Could not find the '...Powered by Jekyll...' footer paragraph.\"\n}\n\n# Delete the blank spacer + \"Ver no Jupiter...\" + \"(c) 2020 ...\" paragraphs\n# in one shot via a Range spanning from just after LOQ4064 through the\n# end of the footer paragraph (its paragraph mark included).\n$startPara = $d.Paragraphs.Item($reqIndex + 1)\n$endPara = $d.Paragraphs.Item($endIndex)\n$range = $d.Range($startPara.Range.Start, $endPara.Range.End)\n$range.Delete()\n"}
